# Update the corrXlsx.xlsx "rawData" sheet with the new data layout:
#   - new leading descriptive columns: congruent, corrAns, letterColor, text
#   - "n" (was column A) moves to column E
#   - the old "index_mean"/"index_std" derived columns are dropped
#   - the old "index_raw" column (C) becomes "order" in column N
#   - the old unlabeled raw column (D) moves to column O
#   - G3 becomes an empty, quote-prefixed cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Saved window position (xWindow/yWindow in the workbook view)
$window = $wb.Windows.Item(1)
$window.Left = 240
$window.Top = 6040

# ---- Header row (row 1) ----
# (Note: "order" (N1) is assigned after the red/green/blue data below so the
#  shared-string table fills up in the same order the original workbook used.)
$ws.Range("A1").Value = "congruent"
$ws.Range("B1").Value = "corrAns"
$ws.Range("C1").Value = "letterColor"
$ws.Range("D1").Value = "text"
$ws.Range("E1").Value = "n"
$ws.Range("F1").Value = "resp_mean"
$ws.Range("G1").Value = "resp_raw"
$ws.Range("I1").Value = "resp_std"
$ws.Range("J1").Value = "rt_mean"
$ws.Range("K1").Value = "rt_raw"
$ws.Range("M1").Value = "rt_std"

# ---- Row 2 ----
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "red"
$ws.Range("D2").Value = "red"
$ws.Range("E2").Value = 2
$ws.Range("N2").Value = 4
$ws.Range("O2").Value = 9

# ---- Row 3 ----
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "green"
$ws.Range("D3").Value = "red"
$ws.Range("E3").Value = 2
$ws.Range("G3").Value = "'"
$ws.Range("G3").Value = ""
$ws.Range("N3").Value = 2
$ws.Range("O3").Value = 10

# ---- Row 4 ----
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "green"
$ws.Range("D4").Value = "green"
$ws.Range("E4").Value = 2
$ws.Range("N4").Value = 3
$ws.Range("O4").Value = 6

# ---- Row 5 ----
$ws.Range("A5").Value = 0
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "blue"
$ws.Range("D5").Value = "green"
$ws.Range("E5").Value = 2
$ws.Range("N5").Value = 5
$ws.Range("O5").Value = 8

# ---- Row 6 ----
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = "blue"
$ws.Range("D6").Value = "blue"
$ws.Range("E6").Value = 2
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 11

# ---- Row 7 ----
$ws.Range("A7").Value = 0
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "red"
$ws.Range("D7").Value = "blue"
$ws.Range("E7").Value = 2
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 7

# "order" header (added to the shared-string table after the color strings,
# matching the original file's string order)
$ws.Range("N1").Value = "order"

# ---- Selection matches the saved state in the target workbook ----
[void]$ws.Range("G3").Select()
